$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 95
$ws.Range("F4").Value = 10176
$ws.Range("F5").Value = 734
$ws.Range("F6").Value = 196
$ws.Range("F7").Value = 417
$ws.Range("F8").Value = 416
$ws.Range("F9").Value = 463
$ws.Range("F11").Value = 12745
$ws.Range("F12").Value = 44
$ws.Range("F18").Value = 46
$ws.Range("F19").Value = 190
$ws.Range("F20").Value = 142
$ws.Range("F22").Value = 2757
$ws.Range("F24").Value = 105
$ws.Range("F27").Value = 69
$ws.Range("F28").Value = 2170
$ws.Range("F29").Value = 1099
$ws.Range("F30").Value = 4261
$ws.Range("F32").Value = 3785
$ws.Range("F33").Value = 835
$ws.Range("F34").Value = 2652
$ws.Range("F35").Value = 3084
$ws.Range("F36").Value = 69
$ws.Range("F37").Value = 1367
$ws.Range("F39").Value = 788
$ws.Range("F43").Value = 668
$ws.Range("F46").Value = 289
$ws.Range("F47").Value = 120
$ws.Range("F48").Value = 158
$ws.Range("F49").Value = 173

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F13").Value = 63

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 95
$ws.Range("F5").Value = 10176
$ws.Range("F6").Value = 734
$ws.Range("F8").Value = 196
$ws.Range("F9").Value = 417
$ws.Range("F10").Value = 416
$ws.Range("F11").Value = 463
$ws.Range("F13").Value = 12745
$ws.Range("F17").Value = 46
$ws.Range("F19").Value = 190
$ws.Range("F20").Value = 142
$ws.Range("F22").Value = 2757
$ws.Range("F24").Value = 105
$ws.Range("F27").Value = 69
$ws.Range("F28").Value = 2170
$ws.Range("F29").Value = 1099
$ws.Range("F30").Value = 4261
$ws.Range("F31").Value = 3785
$ws.Range("F32").Value = 835
$ws.Range("F33").Value = 2652
$ws.Range("F34").Value = 3084
$ws.Range("F35").Value = 69
$ws.Range("F36").Value = 1367
$ws.Range("F38").Value = 788
$ws.Range("F43").Value = 668
$ws.Range("F46").Value = 289
$ws.Range("F47").Value = 120
$ws.Range("F48").Value = 158
$ws.Range("F49").Value = 173
